$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 394, shifting existing rows 394:516 down to 395:517
$ws.Rows.Item(394).Insert()

# Populate the newly inserted row 394 with the new record's data
$ws.Range("A394").Value = 4
$ws.Range("B394").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C394").Value = "Los Lagos"
$ws.Range("D394").Value = 45215
$ws.Range("E394").Value = 10
$ws.Range("F394").Value = 100112043
$ws.Range("G394").Value = "Pepino ensalada"
$ws.Range("H394").Value = "Sin especificar"
$ws.Range("I394").Value = "Primera"
$ws.Range("J394").Value = 80
$ws.Range("K394").Value = 17000
$ws.Range("L394").Value = 17000
$ws.Range("M394").Value = 17000
$ws.Range("N394").Value = "$/caja 60 unidades"
$ws.Range("O394").Value = "Región de Arica y Parinacota"
$ws.Range("P394").Value = 283
$ws.Range("Q394").Value = 60
$ws.Range("R394").Value = "Hortaliza"
